$d = $word.ActiveDocument

$replacements = @(
    @{Old = "Eine Welt des Geschmacks in einer Tasse"; New = "Chai-Tee: Eine Welt des Geschmacks in einer Tasse"},
    @{Old = "Die perfekte Mischung aus Gesundheit und Genuss"; New = "Chai Tee: Die perfekte Mischung aus Gesundheit und Vergnügen"},
    @{Old = "Mehr als nur Tee, eine Lebenseinstellung"; New = "Chai Tee: Mehr als nur Tee, eine Lebensart"},
    @{Old = "Ein Getränk für alle Jahreszeiten und Anlässe"; New = "Chai Tee: Ein Getränk aus allen Jahreszeiten und Gründen"},
    @{Old = "Der ultimative Genuss für Ihre Sinne"; New = "Chai-Tee: Der ultimative Genuss für Ihre Sinne"},
    @{Old = "Eine süße Flucht aus dem Alltag"; New = "Chai-Tee: Eine süße Flucht vom Alltag"},
    @{Old = "Teilen Sie die Wärme, teilen Sie die Liebe"; New = "Chai-Tee: Teilen Sie die Wärme, teilen Sie die Liebe"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
